# Apply edits described by the diff:
# - Row 5 (A5:B5) numeric junk -> "carbon" / "CLIMATE-PROPERTY"
# - Row 12 (A12:B12) "carbon"/"CLIMATE-PROPERTY" -> "atmospheric chemistry"/"not a predefined category"
# - Row 13 (A13:B13) "atmospheric chemistry"/"not a predefined category" -> "Australia"/"CLIMATE-ORGANISATIONS"
# - Row 22 (A22:B22) "Australia"/"CLIMATE-ORGANISATIONS" -> "smallholder farmers"/"CLIMATE-ORGSANISMS"
# - Row 23 deleted entirely (sheet shrinks from A1:D23 to A1:D22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "carbon"
$ws.Range("B5").Value = "CLIMATE-PROPERTY"

$ws.Range("A12").Value = "atmospheric chemistry"
$ws.Range("B12").Value = "not a predefined category"

$ws.Range("A13").Value = "Australia"
$ws.Range("B13").Value = "CLIMATE-ORGANISATIONS"

$ws.Range("A22").Value = "smallholder farmers"
$ws.Range("B22").Value = "CLIMATE-ORGSANISMS"

# Remove the now-duplicate trailing row 23 entirely, shifting rows up.
$ws.Rows(23).Delete()
